$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.838.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.468.62"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.468.56"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.92"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.917.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.804.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.482.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +9.92%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +18.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "656.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.64%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.591.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0979"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.02%  "
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.41"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0315"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -61.28%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.608"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0512"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.46%  "
